$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new job posting was added at the top of the list; insert a row above the
# current row 2 so every existing posting shifts down by one row.
$ws.Rows.Item(2).Insert()
$ws.Rows.Item(2).ClearFormats()

# Populate the newly inserted row 2 with the new posting's details.
$ws.Cells.Item(2, 1).Value = "Oracle Careers"
$ws.Cells.Item(2, 2).Value = "Specialist – Operations, Fund Risk Management and Oversight (Open to all applicants)"
$ws.Cells.Item(2, 3).Value = "New Delhi, India"
$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = "02/04/2026"
$ws.Rows.Item(2).ClearFormats()
$ws.Cells.Item(2, 5).Formula = '=HYPERLINK("https://estm.fa.em2.oraclecloud.com/hcmUI/CandidateExperience/en/sites/CX_1/job/28371/?location=India&locationId=300000000440677&locationLevel=country&mode=location", "Apply")'
